# Generate Report for Handoff
#
# A new handoff was generated for file
# 44e48f6d-14e0-46ea-9fd1-e0dacc693fab, so its "latest handoff"
# timestamps are refreshed on the Overview sheet and on each
# per-language detail sheet (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest Handoff Date" column (D) for the
#     44e48f6d-14e0-46ea-9fd1-e0dacc693fab.md row (row 4) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("D4").Value = "2016-38-18 02:38:02"

# --- zh-cn sheet: "Latest Handoff Datetime" column (E) for the
#     44e48f6d-14e0-46ea-9fd1-e0dacc693fab row (row 4) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-18 02:37:54"

# --- de-de sheet: "Latest Handoff Datetime" column (E) for the
#     44e48f6d-14e0-46ea-9fd1-e0dacc693fab row (row 4) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-18 02:38:02"
